# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Row => new value for sheet "展览" (sheet1)
$exhibitionChanges = @{
    2  = 3
    3  = 167
    5  = 17
    7  = 1649
    8  = 7
    9  = 11
    11 = 1516
    12 = 122
    13 = 45
    14 = 376
    15 = 252
    16 = 190
    18 = 20
    19 = 24
    21 = 269
    22 = 148
    23 = 216
    24 = 208
}

# Row => new value for sheet "全部类型" (sheet4)
$allTypesChanges = @{
    2  = 3
    3  = 167
    5  = 17
    7  = 1649
    9  = 7
    10 = 11
    12 = 1516
    13 = 122
    14 = 45
    15 = 376
    16 = 252
    17 = 190
    19 = 20
    20 = 24
    22 = 269
    23 = 148
    24 = 216
    25 = 208
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionChanges.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionChanges[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesChanges.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesChanges[$row]
}
